# Data update from DGS's 2021/08/25 report.
# Appends a new observation row (row 72) to the risk-matrix time series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 71
$newRow = $lastRow + 1

# Column A holds the report date as literal text (e.g. "2021/08/23"), not a
# real Excel date serial. Temporarily force a text number format so the
# date-shaped string isn't auto-converted to a date serial, then restore the
# same date-display format used by the rest of the column (cosmetic only,
# since the cell's underlying type stays text).
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2021/08/25"
$ws.Range("A" + $newRow).NumberFormat = $ws.Range("A" + $lastRow).NumberFormat

# Columns B:E hold the numeric series values for the new report date.
$ws.Range("B" + $newRow).Value = 312.8
$ws.Range("C" + $newRow).Value = 317.1
$ws.Range("D" + $newRow).Value = 0.98
$ws.Range("E" + $newRow).Value = 0.98

# Mirror the original author's follow-on selection: the cell just below the
# newly appended row.
$ws.Range("A" + ($newRow + 1)).Select()
